$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s = $ws.Range("D2").Style
$ws.Range("D2").Formula = "'244.44"
$ws.Range("D2").Style = $s
$s = $ws.Range("G2").Style
$ws.Range("G2").Formula = "'7"
$ws.Range("G2").Style = $s
$s = $ws.Range("D3").Style
$ws.Range("D3").Formula = "'24.06"
$ws.Range("D3").Style = $s
$s = $ws.Range("G3").Style
$ws.Range("G3").Formula = "'7"
$ws.Range("G3").Style = $s
$s = $ws.Range("D4").Style
$ws.Range("D4").Formula = "'5.209"
$ws.Range("D4").Style = $s
$s = $ws.Range("G4").Style
$ws.Range("G4").Formula = "'7"
$ws.Range("G4").Style = $s
$s = $ws.Range("D5").Style
$ws.Range("D5").Formula = "'0.05776"
$ws.Range("D5").Style = $s
$s = $ws.Range("G5").Style
$ws.Range("G5").Formula = "'7"
$ws.Range("G5").Style = $s
$s = $ws.Range("D6").Style
$ws.Range("D6").Formula = "'6.499"
$ws.Range("D6").Style = $s
$s = $ws.Range("G6").Style
$ws.Range("G6").Formula = "'7"
$ws.Range("G6").Style = $s
$s = $ws.Range("D7").Style
$ws.Range("D7").Formula = "'3.123"
$ws.Range("D7").Style = $s
$s = $ws.Range("G7").Style
$ws.Range("G7").Formula = "'7"
$ws.Range("G7").Style = $s
$s = $ws.Range("D8").Style
$ws.Range("D8").Formula = "'0.8146"
$ws.Range("D8").Style = $s
$s = $ws.Range("G8").Style
$ws.Range("G8").Formula = "'7"
$ws.Range("G8").Style = $s
$s = $ws.Range("D9").Style
$ws.Range("D9").Formula = "'0.8491"
$ws.Range("D9").Style = $s
$s = $ws.Range("G9").Style
$ws.Range("G9").Formula = "'7"
$ws.Range("G9").Style = $s
$s = $ws.Range("D10").Style
$ws.Range("D10").Formula = "'0.1357"
$ws.Range("D10").Style = $s
$s = $ws.Range("G10").Style
$ws.Range("G10").Formula = "'7"
$ws.Range("G10").Style = $s
$s = $ws.Range("D11").Style
$ws.Range("D11").Formula = "'0.06954"
$ws.Range("D11").Style = $s
$s = $ws.Range("G11").Style
$ws.Range("G11").Formula = "'7"
$ws.Range("G11").Style = $s
$s = $ws.Range("D12").Style
$ws.Range("D12").Formula = "'0.03157"
$ws.Range("D12").Style = $s
$s = $ws.Range("G12").Style
$ws.Range("G12").Formula = "'7"
$ws.Range("G12").Style = $s
$s = $ws.Range("D13").Style
$ws.Range("D13").Formula = "'0.02868"
$ws.Range("D13").Style = $s
$s = $ws.Range("G13").Style
$ws.Range("G13").Formula = "'7"
$ws.Range("G13").Style = $s
$s = $ws.Range("D14").Style
$ws.Range("D14").Formula = "'0.09371"
$ws.Range("D14").Style = $s
$s = $ws.Range("G14").Style
$ws.Range("G14").Formula = "'7"
$ws.Range("G14").Style = $s
$s = $ws.Range("G15").Style
$ws.Range("G15").Formula = "'7"
$ws.Range("G15").Style = $s
$s = $ws.Range("D16").Style
$ws.Range("D16").Formula = "'0.001510"
$ws.Range("D16").Style = $s
$s = $ws.Range("G16").Style
$ws.Range("G16").Formula = "'7"
$ws.Range("G16").Style = $s
$s = $ws.Range("D17").Style
$ws.Range("D17").Formula = "'0.04687"
$ws.Range("D17").Style = $s
$s = $ws.Range("G17").Style
$ws.Range("G17").Formula = "'7"
$ws.Range("G17").Style = $s
$s = $ws.Range("D18").Style
$ws.Range("D18").Formula = "'0.0005992"
$ws.Range("D18").Style = $s
$s = $ws.Range("G18").Style
$ws.Range("G18").Formula = "'7"
$ws.Range("G18").Style = $s
$s = $ws.Range("D19").Style
$ws.Range("D19").Formula = "'0.006278"
$ws.Range("D19").Style = $s
$s = $ws.Range("G19").Style
$ws.Range("G19").Formula = "'7"
$ws.Range("G19").Style = $s
$s = $ws.Range("D20").Style
$ws.Range("D20").Formula = "'0.001237"
$ws.Range("D20").Style = $s
$s = $ws.Range("G20").Style
$ws.Range("G20").Formula = "'7"
$ws.Range("G20").Style = $s
$s = $ws.Range("D21").Style
$ws.Range("D21").Formula = "'0.004294"
$ws.Range("D21").Style = $s
$s = $ws.Range("G21").Style
$ws.Range("G21").Formula = "'7"
$ws.Range("G21").Style = $s
$s = $ws.Range("D22").Style
$ws.Range("D22").Formula = "'0.00006505"
$ws.Range("D22").Style = $s
$s = $ws.Range("G22").Style
$ws.Range("G22").Formula = "'7"
$ws.Range("G22").Style = $s
$s = $ws.Range("D23").Style
$ws.Range("D23").Formula = "'3.498"
$ws.Range("D23").Style = $s
$s = $ws.Range("G23").Style
$ws.Range("G23").Formula = "'7"
$ws.Range("G23").Style = $s
$s = $ws.Range("G24").Style
$ws.Range("G24").Formula = "'7"
$ws.Range("G24").Style = $s
$s = $ws.Range("G25").Style
$ws.Range("G25").Formula = "'7"
$ws.Range("G25").Style = $s
$s = $ws.Range("D26").Style
$ws.Range("D26").Formula = "'0.1338"
$ws.Range("D26").Style = $s
$s = $ws.Range("G26").Style
$ws.Range("G26").Formula = "'7"
$ws.Range("G26").Style = $s
$s = $ws.Range("G27").Style
$ws.Range("G27").Formula = "'7"
$ws.Range("G27").Style = $s
$s = $ws.Range("D28").Style
$ws.Range("D28").Formula = "'0.0002330"
$ws.Range("D28").Style = $s
$s = $ws.Range("G28").Style
$ws.Range("G28").Formula = "'7"
$ws.Range("G28").Style = $s
$s = $ws.Range("G29").Style
$ws.Range("G29").Formula = "'7"
$ws.Range("G29").Style = $s
$s = $ws.Range("G30").Style
$ws.Range("G30").Formula = "'7"
$ws.Range("G30").Style = $s
$s = $ws.Range("G31").Style
$ws.Range("G31").Formula = "'7"
$ws.Range("G31").Style = $s
$s = $ws.Range("G32").Style
$ws.Range("G32").Formula = "'7"
$ws.Range("G32").Style = $s
$s = $ws.Range("G33").Style
$ws.Range("G33").Formula = "'7"
$ws.Range("G33").Style = $s
$s = $ws.Range("G34").Style
$ws.Range("G34").Formula = "'7"
$ws.Range("G34").Style = $s
$s = $ws.Range("G35").Style
$ws.Range("G35").Formula = "'7"
$ws.Range("G35").Style = $s
$s = $ws.Range("G36").Style
$ws.Range("G36").Formula = "'7"
$ws.Range("G36").Style = $s
$s = $ws.Range("G37").Style
$ws.Range("G37").Formula = "'7"
$ws.Range("G37").Style = $s
$s = $ws.Range("G38").Style
$ws.Range("G38").Formula = "'7"
$ws.Range("G38").Style = $s
$s = $ws.Range("G39").Style
$ws.Range("G39").Formula = "'7"
$ws.Range("G39").Style = $s
$s = $ws.Range("D40").Style
$ws.Range("D40").Formula = "'0.03632"
$ws.Range("D40").Style = $s
$s = $ws.Range("G40").Style
$ws.Range("G40").Formula = "'7"
$ws.Range("G40").Style = $s
$s = $ws.Range("D41").Style
$ws.Range("D41").Formula = "'0.006294"
$ws.Range("D41").Style = $s
$s = $ws.Range("G41").Style
$ws.Range("G41").Formula = "'7"
$ws.Range("G41").Style = $s
$s = $ws.Range("D42").Style
$ws.Range("D42").Formula = "'0.1050"
$ws.Range("D42").Style = $s
$s = $ws.Range("G42").Style
$ws.Range("G42").Formula = "'7"
$ws.Range("G42").Style = $s
$s = $ws.Range("D43").Style
$ws.Range("D43").Formula = "'0.003202"
$ws.Range("D43").Style = $s
$s = $ws.Range("G43").Style
$ws.Range("G43").Formula = "'7"
$ws.Range("G43").Style = $s
$s = $ws.Range("D44").Style
$ws.Range("D44").Formula = "'0.007512"
$ws.Range("D44").Style = $s
$s = $ws.Range("G44").Style
$ws.Range("G44").Formula = "'7"
$ws.Range("G44").Style = $s
$s = $ws.Range("G45").Style
$ws.Range("G45").Formula = "'7"
$ws.Range("G45").Style = $s
$s = $ws.Range("D46").Style
$ws.Range("D46").Formula = "'0.00000000751"
$ws.Range("D46").Style = $s
$s = $ws.Range("G46").Style
$ws.Range("G46").Formula = "'7"
$ws.Range("G46").Style = $s
$s = $ws.Range("D47").Style
$ws.Range("D47").Formula = "'0.2901"
$ws.Range("D47").Style = $s
$s = $ws.Range("G47").Style
$ws.Range("G47").Formula = "'7"
$ws.Range("G47").Style = $s
$s = $ws.Range("D48").Style
$ws.Range("D48").Formula = "'0.002340"
$ws.Range("D48").Style = $s
$s = $ws.Range("G48").Style
$ws.Range("G48").Formula = "'7"
$ws.Range("G48").Style = $s
$s = $ws.Range("D49").Style
$ws.Range("D49").Formula = "'0.00002101"
$ws.Range("D49").Style = $s
$s = $ws.Range("G49").Style
$ws.Range("G49").Formula = "'7"
$ws.Range("G49").Style = $s
$s = $ws.Range("D50").Style
$ws.Range("D50").Formula = "'0.0002001"
$ws.Range("D50").Style = $s
$s = $ws.Range("G50").Style
$ws.Range("G50").Formula = "'7"
$ws.Range("G50").Style = $s
$s = $ws.Range("G51").Style
$ws.Range("G51").Formula = "'7"
$ws.Range("G51").Style = $s